$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Assignment Status" column (B) for rows 15-21 to "Shared" ---
$ws.Range("B15").Value = "Shared"
$ws.Range("B16").Value = "Shared"
$ws.Range("B17").Value = "Shared"
$ws.Range("B18").Value = "Shared"
$ws.Range("B19").Value = "Shared"
$ws.Range("B20").Value = "Shared"
$ws.Range("B21").Value = "Shared"

# --- 2. Add the newly-shared links in column C (rows 15, 17-21) ---
# Set the cell text to the URL, then attach the hyperlink (Logistic Regression /
# row 16 already has its link, so it is left untouched).
$ws.Range("C15").Value = "https://github.com/dhivyadharani86/Data-Science-Assignments/blob/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/SVM/SVM.ipynb"
$ws.Range("C17").Value = "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/Recommendation%20System"
$ws.Range("C18").Value = "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/XGBM%20%26%20LGBM"
$ws.Range("C19").Value = "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/Neural%20networks"
$ws.Range("C20").Value = "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/NLP%20and%20Naive%20Bayes"
$ws.Range("C21").Value = "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/Timeseries"

$ws.Hyperlinks.Add($ws.Range("C15"), "https://github.com/dhivyadharani86/Data-Science-Assignments/blob/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/SVM/SVM.ipynb")
$ws.Hyperlinks.Add($ws.Range("C17"), "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/Recommendation%20System")
$ws.Hyperlinks.Add($ws.Range("C18"), "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/XGBM%20%26%20LGBM")
$ws.Hyperlinks.Add($ws.Range("C19"), "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/Neural%20networks")
$ws.Hyperlinks.Add($ws.Range("C20"), "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/NLP%20and%20Naive%20Bayes")
$ws.Hyperlinks.Add($ws.Range("C21"), "https://github.com/dhivyadharani86/Data-Science-Assignments/tree/2e364424cec678bf5bf466c2c51702a9d79bddf5/Assignments/Timeseries")

# --- 3. Re-apply the same visual style used by the other hyperlink cells (C16) so ---
#        the newly-linked cells match the rest of the "Link" column formatting.
$ws.Range("C16").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C21").PasteSpecial(-4122)

$excel.CutCopyMode = $false
